$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.476.06"
$ws.Range("E2").Value = "  -1.04%  "
$ws.Range("D3").Value = "2.379.37"
$ws.Range("E3").Value = "  +5.93%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'233.09"
$ws.Range("E5").Value = "  +0.46%  "
$ws.Range("D6").Value = "'0.654"
$ws.Range("E6").Value = "  +1.85%  "
$ws.Range("D7").Value = "'69.58"
$ws.Range("E7").Value = "  +9.82%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "'0.459"
$ws.Range("E9").Value = "  +2.62%  "
$ws.Range("D10").Value = "'0.0966"
$ws.Range("E10").Value = "  -1.35%  "
$ws.Range("D11").Value = "'57.21"
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "2.734.29"
$ws.Range("E12").Value = "  +5.93%  "
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").Value = "'26.22"
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").Value = "'15.65"
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("D16").Value = "'6.21"
$ws.Range("E16").Value = "  +1.49%  "
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").Value = "'0.849"
$ws.Range("E17").Value = "  +2.34%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.373.41"
$ws.Range("E18").Value = "  +5.60%  "
$ws.Range("D19").Value = "43.518.46"
$ws.Range("D20").Value = "0.0₃0985"
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("E21").Value = "  +4.27%  "
$ws.Range("D22").Value = "'74.05"
$ws.Range("E22").Value = "  +1.91%  "
$ws.Range("D23").Value = "'250.88"
$ws.Range("E23").Value = "  +1.26%  "
$ws.Range("D24").Value = "'3.95"
$ws.Range("E24").Value = "  +18.35%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").Value = "  +1.54%  "
$ws.Range("D27").Value = "'2.22"
$ws.Range("E27").Value = "  -3.17%  "
$ws.Range("D28").Value = "'22.99"
$ws.Range("E28").Value = "  +9.45%  "
$ws.Range("D29").Value = "'9.93"
$ws.Range("E29").Value = "  +1.07%  "
$ws.Range("D30").Value = "'172.61"
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("D31").Value = "'1.54"
$ws.Range("E31").Value = "  +9.05%  "
$ws.Range("D32").Value = "'0.126"
$ws.Range("E32").Value = "  -9.06%  "
$ws.Range("E33").Value = "  +1.62%  "
$ws.Range("D34").Value = "'4.95"
$ws.Range("E34").Value = "  +3.38%  "
$ws.Range("E35").Value = "  +0.32%  "
$ws.Range("E36").Value = "  +2.45%  "
$ws.Range("E37").Value = "  +2.84%  "
$ws.Range("D38").Value = "'2.45"
$ws.Range("E38").Value = "  +7.50%  "
$ws.Range("E39").Value = "  -0.85%  "
$ws.Range("D40").Value = "'0.0253"
$ws.Range("E40").Value = "  +0.36%  "
$ws.Range("E41").Value = "  +0.17%  "
$ws.Range("D42").Value = "'8.90"
$ws.Range("E42").Value = "  +4.15%  "
$ws.Range("D43").Value = "'18.39"
$ws.Range("E43").Value = "  +7.69%  "
$ws.Range("E44").Value = "  +10.82%  "
$ws.Range("D45").Value = "'4.52"
$ws.Range("E45").Value = "  +4.96%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").Value = "'1.22"
$ws.Range("E46").Value = "  +1.99%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'99.00"
$ws.Range("E47").Value = "  +1.72%  "
$ws.Range("D48").Value = "'0.0949"
$ws.Range("E48").Value = "  +0.83%  "
$ws.Range("D49").Value = "1.446.55"
$ws.Range("E49").Value = "  +0.58%  "
$ws.Range("D50").Value = "2.605.80"
$ws.Range("E50").Value = "  +6.10%  "
$ws.Range("D51").Value = "'0.000200"
$ws.Range("E51").Value = "  -10.52%  "
